$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date (stored as an Excel
# serial date) for every data row. The whole column (rows 2-510) is bumped
# from 46061 (2026-02-08) to 46062 (2026-02-09).
$ws.Range("C2:C510").Value2 = 46062
